$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Iteration #1")
$ws5 = $wb.Worksheets.Item("Iteration #4")

# Remove one of the currently-empty filler rows (row 36) so the rest of
# the sheet (TOTAL row, auto-evaluation block, ...) shifts up by one to
# make room for the new entry while keeping the same overall layout.
$ws5.Rows.Item(36).Delete()

# Row 19 was an empty filler row; turn it into a real iteration entry.
# Copy the cell formatting used by the other "long task description"
# rows (date cell + wrapped text cell) instead of re-building styles by
# hand, then set the actual values.
$ws1.Range("A19").Copy()
$ws5.Range("A19").PasteSpecial(-4122)
$ws5.Range("A19").Value2 = 42855

$ws1.Range("B19").Copy()
$ws5.Range("B19").PasteSpecial(-4122)
$ws5.Range("B19").Value = "Les examens se terminent selon un chrono + Implémentation d'une animation pour les progress bars."

$ws5.Range("C19").Value = 3

$ws5.Rows.Item(19).RowHeight = 30

$excel.CutCopyMode = 0
